# Refresh the "cryptos" price/volume snapshot (GitHub Actions cron update).
# Price cells (column D) are stored as text even when they look numeric
# (e.g. "236.22", "30.233.40"), so NumberFormat is forced to "@" (Text)
# before writing any value that Excel would otherwise auto-convert to a
# number - this preserves exact formatting (trailing zeros, thousands-style
# dots, etc.) instead of losing it to numeric coercion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.233.40"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.859.18"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.22"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4790"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2803"
$ws.Range("E8").Value = "  -4.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06452"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").Value = "1.859.91"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07389"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.26"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.090"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.11"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6455"
$ws.Range("E15").Value = "  -3.58%  "
$ws.Range("D16").Value = "30.170.59"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9994"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.16"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007568"
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "224.08"
$ws.Range("E20").Value = "  +13.30%  "
$ws.Range("D21").Value = "2.098.77"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.265"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.092"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.212"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.56"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.928"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  -3.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09221"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.235"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.956"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04973"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.144"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7241"
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01834"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.601"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8999"
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.043"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.20"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.880"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4250"
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1305"
$ws.Range("E45").Value = "  -4.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.275"
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "63.78"
$ws.Range("E47").Value = "  -8.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.498"
$ws.Range("E48").Value = "  +6.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.685"
$ws.Range("E49").Value = "  -3.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.75"
$ws.Range("E50").Value = "  -3.90%  "
$ws.Range("E51").Value = "  -3.43%  "
